$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.011.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.641.21'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5049'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.009'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2573'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06424'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07722'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.647.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.249'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.869.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5447'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7899'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.039.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.009'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.294'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.965'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.010'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.934'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1154'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.738'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05051'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.242'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.249'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.188'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.542'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.342'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8947'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.613'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5627'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.147.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01571'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.568'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.009'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.679'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8129'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.10%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.779.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈113'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4532'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05041'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.01%  '
